$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Capture the bold "Meta description" lead-in run (empty run + bold run)
#    immediately, since it is the only bold run in the document and we will
#    reuse its structure/formatting to author the new closing paragraph
#    further down ("Play Diamond Multiplier Respin Free - Review & Ratings").
#    This must happen before any other edit touches the document, because a
#    captured FormattedText tracks its live position in the document.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$templateSrc = $metaPara.Range.Duplicate
$templateSrc.End = $templateSrc.Start + 16
$template = $templateSrc.FormattedText

# ---------------------------------------------------------------------------
# 2. Insert and populate the new closing paragraph right away (still using
#    the original, pre-edit paragraph index for the "What we don't like"
#    list's last bullet, "No Autoplay function available") so the template
#    capture above stays valid.
# ---------------------------------------------------------------------------
$bulletPara = $d.Paragraphs(50)
$bulletPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(51)
$newPara.Style = "Normal"
$dest = $newPara.Range.Duplicate
$dest.Collapse(1)
$dest.FormattedText = $template
$d.Paragraphs(51).Range.Find.Execute(
    "Meta description", $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Diamond Multiplier Respin Free - Review & Ratings", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Remove the "Meta description" paragraph entirely.
# ---------------------------------------------------------------------------
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 4. Update the title heading.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Play Diamond Multiplier Respin for Free - Review", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Play Diamond Multiplier Respin Free - Review & Ratings", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Update the "What we like" bullet list.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Unique Wild symbol with multipliers up to 5x", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Simple gameplay with a 3x3 grid and 9 fixed paylines", 2) | Out-Null

$d.Content.Find.Execute(
    "Triggered Respin feature on Wild winning combinations", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Crisp and high definition graphics with vibrant colors", 2) | Out-Null

$d.Content.Find.Execute(
    "Scatter offers 10 free spins with active Wild Respin feature ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Wild symbols with multipliers trigger respins", 2) | Out-Null

$d.Content.Find.Execute(
    "Compatible with desktop and mobile devices thanks to HTML5 technology", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Available for both desktop and mobile devices", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. Update the "What we don't like" bullet list.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Only 5 base symbols", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "No Autoplay function", 2) | Out-Null

$d.Content.Find.Execute(
    "No Autoplay function available", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Limited special features", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7. Replace the DALLE prompt paragraph with the new closing summary text.
# ---------------------------------------------------------------------------
$dalleOld = "Prompt for DALLE: Create a feature image for the game " + [char]34 + `
    "Diamond Multiplier Respin" + [char]34 + " in a cartoon style with a happy " + `
    "Maya warrior wearing glasses. The warrior should be holding a diamond and " + `
    "surrounded by the game's symbols, including the red and golden number 7s, " + `
    "the three different BAR signs, and the Free Spins and Wild symbols. Use " + `
    "vibrant colors and make the image lively to represent the game's potential " + `
    "for high payouts and the energetic '80s dance tune."

$d.Content.Find.Execute(
    $dalleOld, $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Diamond Multiplier Respin and play for free. Find out about the gameplay, graphics, and special features.",
    2) | Out-Null

Write-Output "Done"
